# March 24 update 3
# Adds three new columns (M: renewd, N: PlanID, O: iteration) to Sheet1,
# populating the header row and all data rows (2-13) with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header values
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Match the header formatting (bold, bordered, centered style) used by the
# other header cells (e.g. L1) by copying its format onto the new headers.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for each of the 12 data rows
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"
    $ws.Cells.Item($r, 14).Value = 20141086
    $ws.Cells.Item($r, 15).Value = 2
}
